$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty Average/Max formula cells for rows 4-8 (D:G),
# and for row 6-8 (B:C) which were also blank.

$ws.Range("D4").Formula = "=AVERAGE(H11:H45)"
$ws.Range("E4").Formula = "=MAX(H11:H45)"
$ws.Range("F4").Formula = "=AVERAGE(B37:B86)"
$ws.Range("G4").Formula = "=MAX(B37:B86)"

$ws.Range("D5").Formula = "=AVERAGE(I11:I45)"
$ws.Range("E5").Formula = "=MAX(I11:I45)"
$ws.Range("F5").Formula = "=AVERAGE(C37:C86)"
$ws.Range("G5").Formula = "=MAX(C37:C86)"

$ws.Range("B6").Formula = "=AVERAGE(D11:D35)"
$ws.Range("C6").Formula = "=MAX(D11:D35)"
$ws.Range("D6").Formula = "=AVERAGE(J11:J45)"
$ws.Range("E6").Formula = "=MAX(J11:J45)"
$ws.Range("F6").Formula = "=AVERAGE(D37:D86)"
$ws.Range("G6").Formula = "=MAX(D37:D86)"

$ws.Range("B7").Formula = "=AVERAGE(E11:E35)"
$ws.Range("C7").Formula = "=MAX(E11:E35)"
$ws.Range("D7").Formula = "=AVERAGE(K11:K45)"
$ws.Range("E7").Formula = "=MAX(K11:K45)"
$ws.Range("F7").Formula = "=AVERAGE(E37:E86)"
$ws.Range("G7").Formula = "=MAX(E37:E86)"

$ws.Range("B8").Formula = "=AVERAGE(F11:F35)"
$ws.Range("C8").Formula = "=MAX(F11:F35)"
$ws.Range("D8").Formula = "=AVERAGE(L11:L45)"
$ws.Range("E8").Formula = "=MAX(L11:L45)"
$ws.Range("F8").Formula = "=AVERAGE(F37:F86)"
$ws.Range("G8").Formula = "=MAX(F37:F86)"

# D5:E8 switch to the "box top/bottom-less" border style (matches D4:G4),
# which for this sheet is achieved by copying the border style used in D4.
$ws.Range("D4:E4").Copy() | Out-Null
$ws.Range("D5:E8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update the selected cell to match the final state.
$ws.Range("G13").Select() | Out-Null
